# Update data for carjacking-by-neighborhood-by-month.xlsx
# - Rename the "as of" date from 2022-03-11 to 2022-03-12 (sheet tab + header cell)
# - Apply assorted count corrections/additions across historical columns
# - Insert a new "Galewood" neighborhood row (alphabetically between "Gage Park"
#   and "Garfield Ridge") with its March-2022-to-date count

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet tab and refresh the "through" date in the header row ---
$ws.Name = "Through 2022-03-12"
$ws.Range("B1").Value = "March 2022 (through March 12)"

# --- Numeric corrections/additions (row numbers are pre-insertion positions) ---
$ws.Range("B3").Value = 5     # Austin, March 2022 (through March 12)
$ws.Range("T4").Value = 1     # North Lawndale, March 2016
$ws.Range("B5").Value = 5     # Garfield Park, March 2022 (through March 12)
$ws.Range("Q5").Value = 5     # Garfield Park, March 2017
$ws.Range("H7").Value = 1     # South Shore, March 2020
$ws.Range("Q8").Value = 1     # Kenwood, March 2017
$ws.Range("N9").Value = 2     # Chicago Lawn, March 2018
$ws.Range("E12").Value = 1    # Lake View, March 2021
$ws.Range("Q15").Value = 3    # Humboldt Park, March 2017
$ws.Range("B18").Value = 2    # Washington Heights, March 2022 (through March 12)
$ws.Range("N22").Value = 1    # Chatham, March 2018
$ws.Range("H34").Value = 1    # River North, March 2020
$ws.Range("B36").Value = 1    # Roseland, March 2022 (through March 12)
$ws.Range("B39").Value = 1    # North Center, March 2022 (through March 12)
$ws.Range("H43").Value = 1    # Hermosa, March 2020
$ws.Range("E84").Value = 3    # South Deering, March 2021

# --- Insert the new "Galewood" row (between Gage Park [62] and Garfield Ridge [63]) ---
$ws.Rows.Item(63).Insert()
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A63").Value = "Galewood"
$ws.Range("B63").Value = 1
